# Workbook already open
$wb = $excel.ActiveWorkbook

# --- 1. Rename header on "Weekly Quantity" sheet (sheet 1) ---
$wsWeekly = $wb.Worksheets.Item(1)
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# --- 2. Rename header on "Monthly Trend" sheet (sheet 2) ---
$wsMonthly = $wb.Worksheets.Item(2)
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3. Add a new "PO Forecast" worksheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$headerRange = $wsForecast.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data rows: ds, PO_Forecast, yhat_lower, yhat_upper
$data = @(
    @(44948.99999999999, 11, 2.571175959653114, 19.19788143352253),
    @(44955.99999999999, 11, 3.063885435119403, 19.22946709276879),
    @(44962.99999999999, 11, 3.266379212234483, 18.87575145318853),
    @(44983.99999999999, 11, 2.663898564489628, 18.56553673924199),
    @(44990.99999999999, 11, 2.961265168892732, 18.94044119077631),
    @(45025.99999999999, 12, 3.825520547896382, 19.48108610005115),
    @(45053.99999999999, 12, 3.992139428554501, 19.79802279937447),
    @(45074.99999999999, 12, 4.635290970560285, 20.07130286149998),
    @(45088.99999999999, 13, 4.829850550088094, 20.23088180400762),
    @(45095.99999999999, 13, 4.504873198166251, 20.54734226975554),
    @(45102.99999999999, 13, 4.968336072421883, 20.63181142853068),
    @(45109.99999999999, 13, 4.926694683586494, 20.20053958785995),
    @(45116.99999999999, 13, 5.250084115514971, 20.60501053306235),
    @(45123.99999999999, 13, 5.105349949577461, 21.12255772826849),
    @(45130.99999999999, 13, 4.796681876311222, 21.08436836865985),
    @(45137.99999999999, 13, 5.783794630916547, 21.32321548069335),
    @(45144.99999999999, 13, 5.213239605011161, 20.68958192924975),
    @(45151.99999999999, 13, 5.686617005750166, 21.49879352589383),
    @(45158.99999999999, 14, 5.315603455766324, 21.62103766796849)
)

$row = 2
foreach ($item in $data) {
    $wsForecast.Cells.Item($row, 1).Value = $item[0]
    $wsForecast.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $wsForecast.Cells.Item($row, 2).Value = $item[1]
    $wsForecast.Cells.Item($row, 3).Value = $item[2]
    $wsForecast.Cells.Item($row, 4).Value = $item[3]
    $row++
}

# Activate the first sheet to mirror the original workbook's selected tab
$wsWeekly.Activate()
